$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Handle extraneous empty cells: collapse the per-line subtitle example
# values (Line_1 / Line_2 / Line_3) into a single shared placeholder value.
$ws.Range("D2").Value = "Example"
$ws.Range("E2").Value = "Example"
$ws.Range("F2").Value = "Example"
